$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 381, pushing the existing rows 381-460 down to 382-461
$ws.Rows("381:381").Insert()

# Populate the newly inserted row 381 with its data
$ws.Range("A381").Value = 4
$ws.Range("B381").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C381").Value = "Los Lagos"
$ws.Range("D381").Value = 45015
$ws.Range("E381").Value = 10
$ws.Range("F381").Value = "Fruta"
$ws.Range("G381").Value = 100104
$ws.Range("H381").Value = "Frutos de pepita"
$ws.Range("I381").Value = 100104005
$ws.Range("J381").Value = "Pera"
$ws.Range("K381").Value = "Forelle"
$ws.Range("L381").Value = "Primera"
$ws.Range("M381").Value = 400
$ws.Range("N381").Value = 13000
$ws.Range("O381").Value = 14000
$ws.Range("P381").Value = 13500
$ws.Range("Q381").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R381").Value = "Región de O'Higgins"
$ws.Range("S381").Value = 900
$ws.Range("T381").Value = 15
